$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-23 23:18:20"
$ws.Range("O2").Value = "5.5 °C"
$ws.Range("E3").Value = "2026-02-23 23:18:22"
$ws.Range("E4").Value = "2026-02-23 23:18:25"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "71%"
$ws.Range("O4").Value = "11.8 °C"
$ws.Range("E5").Value = "2026-02-23 23:18:27"
$ws.Range("E6").Value = "2026-02-23 23:18:29"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "63%"
$ws.Range("J6").Value = "1024.3 hPa"
$ws.Range("O6").Value = "13.8 °C"
$ws.Range("E7").Value = "2026-02-23 23:18:32"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "68%"
$ws.Range("J7").Value = "1024.6 hPa"
$ws.Range("E8").Value = "2026-02-23 23:18:34"
$ws.Range("J8").Value = "1024.2 hPa"
$ws.Range("E9").Value = "2026-02-23 23:18:37"
$ws.Range("O9").Value = "12.2 °C"
$ws.Range("E10").Value = "2026-02-23 23:18:39"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "78%"
$ws.Range("O10").Value = "10.4 °C"
$ws.Range("E11").Value = "2026-02-23 23:18:42"
$ws.Range("O11").Value = "8.5 °C"
$ws.Range("E12").Value = "2026-02-23 23:18:43"
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-23 23:18:46"
$ws.Range("J13").Value = "1027.0 hPa"
$ws.Range("O13").Value = "6.8 °C"
$ws.Range("E14").Value = "2026-02-23 23:18:48"
$ws.Range("O14").Value = "12.1 °C"
$ws.Range("E15").Value = "2026-02-23 23:18:51"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "72%"
$ws.Range("O15").Value = "12.1 °C"
$ws.Range("E16").Value = "2026-02-23 23:18:53"
$ws.Range("E17").Value = "2026-02-23 23:18:55"
$ws.Range("E18").Value = "2026-02-23 23:18:58"
$ws.Range("O18").Value = "10.5 °C"
$ws.Range("E19").Value = "2026-02-23 23:19:00"
$ws.Range("E20").Value = "2026-02-23 23:19:03"
$ws.Range("O20").Value = "3.9 °C"
$ws.Range("E21").Value = "2026-02-23 23:19:05"
$ws.Range("O21").Value = "9.5 °C"
$ws.Range("E22").Value = "2026-02-23 23:19:07"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "25%"
$ws.Range("E23").Value = "2026-02-23 23:19:10"
$ws.Range("E24").Value = "2026-02-23 23:19:12"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "84%"
$ws.Range("J24").Value = "1025.9 hPa"
$ws.Range("E25").Value = "2026-02-23 23:19:15"
$ws.Range("O25").Value = "5.7 °C"
$ws.Range("E26").Value = "2026-02-23 23:19:17"
$ws.Range("O26").Value = "9.8 °C"
$ws.Range("E27").Value = "2026-02-23 23:19:20"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "30%"
$ws.Range("O27").Value = "5.6 °C"
$ws.Range("E28").Value = "2026-02-23 23:19:22"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "70%"
$ws.Range("O28").Value = "10.8 °C"
$ws.Range("E29").Value = "2026-02-23 23:19:24"
$ws.Range("O29").Value = "10.4 °C"
$ws.Range("E30").Value = "2026-02-23 23:19:27"
$ws.Range("O30").Value = "12.7 °C"
$ws.Range("E31").Value = "2026-02-23 23:19:29"
$ws.Range("E32").Value = "2026-02-23 23:19:32"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "70%"
$ws.Range("N32").Value = "-2.0 °C 22:57 TU"
$ws.Range("O32").Value = "6.8 °C"
$ws.Range("E33").Value = "2026-02-23 23:19:34"
$ws.Range("O33").Value = "8.5 °C"
$ws.Range("E34").Value = "2026-02-23 23:19:37"
$ws.Range("E35").Value = "2026-02-23 23:19:39"
$ws.Range("O35").Value = "11.8 °C"
$ws.Range("E36").Value = "2026-02-23 23:19:41"
$ws.Range("O36").Value = "12.8 °C"
$ws.Range("E37").Value = "2026-02-23 23:19:44"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "69%"
$ws.Range("O37").Value = "8.7 °C"
$ws.Range("E38").Value = "2026-02-23 23:19:46"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "66%"
$ws.Range("O38").Value = "12.0 °C"
$ws.Range("E39").Value = "2026-02-23 23:19:48"
$ws.Range("K39").Value = "16.6 MJ/m2"
$ws.Range("E40").Value = "2026-02-23 23:19:51"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "63%"
$ws.Range("O40").Value = "8.4 °C"
$ws.Range("E41").Value = "2026-02-23 23:19:53"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "75%"
$ws.Range("O41").Value = "11.6 °C"
$ws.Range("E42").Value = "2026-02-23 23:19:55"
$ws.Range("O42").Value = "11.6 °C"
$ws.Range("E43").Value = "2026-02-23 23:19:57"
$ws.Range("O43").Value = "10.3 °C"
$ws.Range("E44").Value = "2026-02-23 23:20:00"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "36%"
$ws.Range("N44").Value = "-0.3 °C 22:59 TU"
$ws.Range("E45").Value = "2026-02-23 23:20:02"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "52%"
$ws.Range("O45").Value = "8.0 °C"
$ws.Range("E46").Value = "2026-02-23 23:20:05"
$ws.Range("O46").Value = "9.9 °C"
